# Cover slide (slide 1): remove the decorative "dot" ovals from the
# gradient panel and tidy up the remaining shapes' display names so
# their numbering is contiguous again (matches how PowerPoint
# renames shapes after the gap left by the deleted ovals).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Remove the four decorative dot shapes -----------------------
# Walk backwards so deleting doesn't shift the indices of shapes we
# still need to visit.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Oval 2" -or $shape.Name -eq "Oval 3" -or `
        $shape.Name -eq "Oval 4" -or $shape.Name -eq "Oval 5") {
        $shape.Delete()
    }
}

# --- 2. Rename the remaining shapes ----------------------------------
$renames = @{
    "Picture 6"            = "Picture 2"
    "Picture 7"            = "Picture 3"
    "TextBox 8"            = "TextBox 4"
    "TextBox 9"            = "TextBox 5"
    "Rounded Rectangle 10" = "Rounded Rectangle 6"
    "TextBox 11"           = "TextBox 7"
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($renames.ContainsKey($shape.Name)) {
        $shape.Name = $renames[$shape.Name]
    }
}
